$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "35.445.15"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "1.903.09"
$ws.Range("E3").Value = "  +2.94%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.64%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.45"
$ws.Range("E5").Value = "  +2.86%  "

# Row 6: 'XRP' -> 'XRP'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.650"
$ws.Range("E6").Value = "  +5.03%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = "  +0.70%  "

# Row 8: 'Solana' -> 'Solana'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.56"
$ws.Range("E8").Value = "  -0.98%  "

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.344"
$ws.Range("E9").Value = "  +5.16%  "

# Row 10: 'OKB' -> 'OKB'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.83"
$ws.Range("E10").Value = "  +7.41%  "

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0711"
$ws.Range("E11").Value = "  +2.88%  "

# Row 12: 'TRON' -> 'TRON'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0999"
$ws.Range("E12").Value = "  +1.08%  "

# Row 13: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D13").Value = "2.184.26"
$ws.Range("E13").Value = "  +3.07%  "

# Row 14: 'Chainlink' -> 'Chainlink'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.14"
$ws.Range("E14").Value = "  +7.00%  "

# Row 15: 'Polygon' -> 'WrappedEther'
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.940.39"
$ws.Range("E15").Value = "  +3.92%  "

# Row 16: 'Polkadot' -> 'Polygon'
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.694"
$ws.Range("E16").Value = "  +3.20%  "

# Row 17: 'WrappedEther' -> 'Polkadot'
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.85"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D18").Value = "35.398.49"
$ws.Range("E18").Value = "  +1.12%  "

# Row 19: 'Litecoin' -> 'Litecoin'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.83"
$ws.Range("E19").Value = "  +2.68%  "

# Row 20: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D20").Value = "0.0₃0815"
$ws.Range("E20").Value = "  +2.98%  "

# Row 21: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.97"
$ws.Range("E21").Value = "  +0.82%  "

# Row 22: 'Avalanche' -> 'Avalanche'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.53"
$ws.Range("E22").Value = "  +3.43%  "

# Row 23: 'Uniswap' -> 'Uniswap'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.76"
$ws.Range("E23").Value = "  +0.30%  "

# Row 24: 'Dai' -> 'Dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.02"
$ws.Range("E24").Value = "  +0.76%  "

# Row 25: 'Toncoin' -> 'Toncoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  +0.72%  "

# Row 26: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +23.01%  "

# Row 27: 'Monero' -> 'Monero'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.48"
$ws.Range("E27").Value = "  +0.90%  "

# Row 28: 'Cosmos' -> 'Cosmos'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.37"
$ws.Range("E28").Value = "  +4.76%  "

# Row 29: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.19"
$ws.Range("E29").Value = "  +3.62%  "

# Row 30: 'Stellar' -> 'Stellar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  +2.10%  "

# Row 31: 'Filecoin' -> 'Filecoin'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.12"
$ws.Range("E31").Value = "  +3.35%  "

# Row 32: 'Hedera' -> 'Hedera'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0564"
$ws.Range("E32").Value = "  +1.66%  "

# Row 33: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  +0.58%  "

# Row 34: 'ImmutableX' -> 'InternetComputer(DFINITY)'
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.14"
$ws.Range("E34").Value = "  +3.52%  "

# Row 35: 'InternetComputer(DFINITY)' -> 'ImmutableX'
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.923"
$ws.Range("E35").Value = "  +19.00%  "

# Row 36: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.75"
$ws.Range("E36").Value = "  +5.36%  "

# Row 37: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("E37").Value = "  +2.81%  "

# Row 38: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.32"
$ws.Range("E38").Value = "  +1.35%  "

# Row 39: 'VeChain' -> 'VeChain'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0210"
$ws.Range("E39").Value = "  +4.77%  "

# Row 40: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.09"
$ws.Range("E40").Value = "  +1.72%  "

# Row 41: 'Kaspa' -> 'Kaspa'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0634"
$ws.Range("E41").Value = "  +14.27%  "

# Row 42: 'Aave' -> 'Aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "89.87"
$ws.Range("E42").Value = "  -0.18%  "

# Row 43: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.71"
$ws.Range("E43").Value = "  +6.76%  "

# Row 44: 'Maker' -> 'Maker'
$ws.Range("D44").Value = "1.340.46"
$ws.Range("E44").Value = "  -0.16%  "

# Row 45: 'RenderToken' -> 'RenderToken'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  +2.23%  "

# Row 46: 'MultiversX' -> 'MultiversX'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.43"
$ws.Range("E46").Value = "  +38.31%  "

# Row 47: 'Gas' -> 'HuobiToken'
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +0.10%  "

# Row 48: 'MXToken' -> 'MXToken'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  +1.74%  "

# Row 49: 'HuobiToken' -> 'Gas'
$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.43"
$ws.Range("E49").Value = "  -7.52%  "

# Row 50: 'FraxShare' -> 'FraxShare'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.53"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51: 'RocketPoolETH' -> 'RocketPoolETH'
$ws.Range("D51").Value = "2.088.52"
$ws.Range("E51").Value = "  +2.14%  "
